# Add a new "City" column (K) to the real-estate data sheet, filling every
# data row with "San Luis Obispo".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("K1").Value = "City"

# Data rows 2-39
$ws.Range("K2:K39").Value = "San Luis Obispo"

# Reflect the scrolled/selected state from the source edit (best effort).
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("M48").Select()
